# Update "想去人数" (F column) figures across the three sheets that
# contain event listings (展览, 演出, 全部类型). 本地生活 is untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    4  = 5319
    6  = 220
    8  = 8924
    10 = 645
    12 = 2611
    13 = 2611
    14 = 6346
    15 = 2342
    17 = 12
    20 = 26
    22 = 6600
    23 = 218
    24 = 82
    25 = 155
    28 = 7238
    31 = 240
    32 = 43
    35 = 28
    40 = 2554
    43 = 14
    44 = 1136
    46 = 558
    47 = 3570
    48 = 104
    49 = 1140
    50 = 34
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2  = 23
    5  = 213
    7  = 96
    15 = 161
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 5319
    4  = 5319
    6  = 220
    7  = 8924
    9  = 645
    10 = 23
    11 = 2611
    14 = 213
    15 = 6346
    16 = 96
    18 = 12
    21 = 26
    24 = 6600
    25 = 218
    27 = 82
    28 = 155
    31 = 7238
    34 = 43
    43 = 14
    44 = 1136
    46 = 3570
    47 = 104
    49 = 1140
    51 = 34
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
